$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old row 11 ("a") - now replaced/merged into new content structure
$ws.Rows(11).Delete()

# Update C6 text
$ws.Range("C6").Value = "File.separator; because the separator is differernt between Windows and Linux;"

# Row 7: page loading performance / performance.timing rich text
$ws.Range("B7").Value = "page loading performance"
$run1 = "Executing performance.timing, provides lots of data that can be used to measure the performance of website:`n"
$run2 = "Page Load: Time needed to load the page (from navigationStart to loadEventEnd).`nWaiting: Waiting for response time (from requestStart to responseStart)`nReceiving: Time needed to download the response (from responseStart to responseEnd)`nDOM Processing: Time needed to build DOM (from responseEnd to loadEventStart`nDOMContentLoaded: Time needed to handle DOMContentLoaded event (from domContentLoadedEventStart to domContentLoadedEventEnd)`nonLoad: Time needed to handle onLoad event (from loadEventStart to loadEventEnd)"
$ws.Range("C7").Value = ($run1 + $run2)
$ws.Range("C7").WrapText = $true
$chars = $ws.Range("C7").Characters($run1.Length + 1, $run2.Length)
$chars.Font.Color = 12419407
$chars.Font.Name = "宋体"
$ws.Rows(7).RowHeight = 121.5

# Row 8: static block
$ws.Range("B8").Value = "static bloclk"
$ws.Range("B8").WrapText = $true
$ws.Range("C8").Value = "Execute immediately when a class been loaded. And the static block only can be execute one time;"

# Row 9: anonymous class
$ws.Range("B9").Value = "new interface(){}"
$ws.Range("C9").Value = "new one anonymous class to implement the interface."

# Row 10: drag element
$ws.Range("B10").Value = "drag element"
$ws.Range("B10").WrapText = $true
$ws.Range("C10").Value = "new Action(driver).dragAndDropBy(draggable,0,10).build.perform;"

# Update selection / view
$ws.Range("C18").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
